# Update market-price / profit figures in Sheets (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

# ALC!row6
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 165.66667
$ws.Range("I6").Value = 165.66667
$ws.Range("K6").Value = 497.00001
$ws.Range("M6").Value = -385.00001

# ALC!row9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 26
$ws.Range("I9").Value = 24.5
$ws.Range("K9").Value = 24.5
$ws.Range("M9").Value = 144.5

# ALC!row33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 310
$ws.Range("I33").Value = 311.75
$ws.Range("K33").Value = 311.75
$ws.Range("M33").Value = -82.75

# ALC!row38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 5118.1
$ws.Range("J38").Value = 16666.334
$ws.Range("L38").Value = 49999.00199999999
$ws.Range("N38").Value = -50743.00199999999

# ALC!row39
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 131.6
$ws.Range("I39").Value = 131.6
$ws.Range("K39").Value = 394.8
$ws.Range("M39").Value = -98.79999999999995

# ALC!row74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5250
$ws.Range("I74").Value = 5250
$ws.Range("K74").Value = 5250
$ws.Range("M74").Value = -4314

# ALC!row77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 5250
$ws.Range("I77").Value = 5250
$ws.Range("K77").Value = 26250
$ws.Range("M77").Value = -21570

# ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3780.1333
$ws.Range("I137").Value = 1911.5
$ws.Range("J137").Value = 5915.7144
$ws.Range("K137").Value = 5734.5
$ws.Range("L137").Value = 17747.1432
$ws.Range("M137").Value = -3184.5
$ws.Range("N137").Value = -22847.1432

# ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4101.522
$ws.Range("I138").Value = 2333
$ws.Range("J138").Value = 4224.9067
$ws.Range("K138").Value = 6999
$ws.Range("L138").Value = 12674.7201
$ws.Range("M138").Value = -1859
$ws.Range("N138").Value = -22954.7201

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6618.2954
$ws.Range("I32").Value = 4780.125
$ws.Range("K32").Value = 4780.125
$ws.Range("M32").Value = -4493.125

# ARM!row38
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()

# ARM!row95
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 50498.5
$ws.Range("J95").Value = 50498.5
$ws.Range("L95").Value = 50498.5
$ws.Range("N95").Value = -55990.5

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2409.7778
$ws.Range("I132").Value = 2089.8696
$ws.Range("J132").Value = 4249.25
$ws.Range("K132").Value = 6269.6088
$ws.Range("L132").Value = 12747.75
$ws.Range("M132").Value = -3739.6088
$ws.Range("N132").Value = -17807.75

# BSM!row7
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 66666700
$ws.Range("I7").Value = 66666700
$ws.Range("K7").Value = 66666700
$ws.Range("M7").Value = -66666587

# BSM!row23
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 900
$ws.Range("I23").Value = 1000
$ws.Range("J23").Value = 800
$ws.Range("K23").Value = 1000
$ws.Range("L23").Value = 800
$ws.Range("M23").Value = -717
$ws.Range("N23").Value = -1366

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1517.2
$ws.Range("I134").Value = 1126.7826
$ws.Range("K134").Value = 3380.3478
$ws.Range("M134").Value = -845.3478

# CRP!row35
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 1274.8572
$ws.Range("I35").Value = 1320.6666
$ws.Range("J35").Value = 1000
$ws.Range("K35").Value = 1320.6666
$ws.Range("L35").Value = 1000
$ws.Range("M35").Value = -1026.6666
$ws.Range("N35").Value = -1588

# CRP!row62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 81960.8
$ws.Range("I62").Value = 2435
$ws.Range("J62").Value = 201249.5
$ws.Range("K62").Value = 2435
$ws.Range("L62").Value = 201249.5
$ws.Range("M62").Value = -1811
$ws.Range("N62").Value = -202497.5

# CRP!row65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 81960.8
$ws.Range("I65").Value = 2435
$ws.Range("J65").Value = 201249.5
$ws.Range("K65").Value = 12175
$ws.Range("L65").Value = 1006247.5
$ws.Range("M65").Value = -9055
$ws.Range("N65").Value = -1012487.5

# CRP!row69
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 22265.428
$ws.Range("I69").Value = 22265.428
$ws.Range("K69").Value = 22265.428
$ws.Range("M69").Value = -21516.428

# CRP!row72
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H72").Value = 22265.428
$ws.Range("I72").Value = 22265.428
$ws.Range("K72").Value = 66796.284
$ws.Range("M72").Value = -63052.284

# CRP!row107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 564.7917
$ws.Range("I107").Value = 370.2353
$ws.Range("J107").Value = 1037.2858
$ws.Range("K107").Value = 370.2353
$ws.Range("L107").Value = 1037.2858
$ws.Range("M107").Value = 1549.7647
$ws.Range("N107").Value = -4877.2858

# CRP!row133
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 66596
$ws.Range("I133").Value = 41196
$ws.Range("K133").Value = 41196
$ws.Range("M133").Value = -38666

# CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2549.1853
$ws.Range("I134").Value = 2059.9092
$ws.Range("K134").Value = 6179.7276
$ws.Range("M134").Value = -3644.7276

# CUL!row132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1673.75
$ws.Range("I132").Value = 1731.6666
$ws.Range("K132").Value = 15584.9994
$ws.Range("M132").Value = -13054.9994

# GSM!row92
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 15875
$ws.Range("J92").Value = 14050
$ws.Range("L92").Value = 14050
$ws.Range("N92").Value = -17794

# LTW!row16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2671.75
$ws.Range("I16").Value = 1612.3334
$ws.Range("K16").Value = 1612.3334
$ws.Range("M16").Value = -1442.3334

# LTW!row31
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 1250
$ws.Range("I31").Value = 166.66667
$ws.Range("J31").Value = 2875
$ws.Range("K31").Value = 166.66667
$ws.Range("L31").Value = 2875
$ws.Range("M31").Value = 81.33332999999999
$ws.Range("N31").Value = -3371

# LTW!row32
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 2037.6666
$ws.Range("I32").Value = 1056.5
$ws.Range("J32").Value = 4000
$ws.Range("K32").Value = 1056.5
$ws.Range("L32").Value = 4000
$ws.Range("M32").Value = -739.5
$ws.Range("N32").Value = -4634

# LTW!row46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1250
$ws.Range("I46").Value = 1250
$ws.Range("K46").Value = 1250
$ws.Range("M46").Value = -1062

# WVR!row10
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 1000000
$ws.Range("I10").Value = 1000000
$ws.Range("K10").Value = 1000000
$ws.Range("M10").Value = -999831

# WVR!row81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9070.909
$ws.Range("J81").Value = 8796.833000000001
$ws.Range("L81").Value = 17593.666
$ws.Range("N81").Value = -19715.666

# WVR!row84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 9070.909
$ws.Range("J84").Value = 8796.833000000001
$ws.Range("L84").Value = 87968.33
$ws.Range("N84").Value = -98576.33
